$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.047.32"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "3.053.81"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.33"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.02"
$ws.Range("E6").Value = "  -1.59%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("E8").Value = "  -1.92%  "

$ws.Range("D9").Value = "3.055.96"
$ws.Range("E9").Value = "  -0.56%  "

$ws.Range("E10").Value = "  -2.84%  "

$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("E12").Value = "  -2.82%  "

$ws.Range("E13").Value = "  -2.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.25"
$ws.Range("E14").Value = "  -3.76%  "

$ws.Range("E15").Value = "  +1.75%  "

$ws.Range("D16").Value = "3.559.72"
$ws.Range("E16").Value = "  -0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.13"
$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("D18").Value = "63.097.89"
$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("D19").Value = "3.055.10"
$ws.Range("E19").Value = "  -0.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.57"
$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("E21").Value = "  -2.71%  "

$ws.Range("E22").Value = "  -1.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.51"
$ws.Range("E23").Value = "  -0.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  -0.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.05"
$ws.Range("E25").Value = "  +1.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.66"
$ws.Range("E26").Value = "  -2.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.49"
$ws.Range("E27").Value = "  +3.78%  "

$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.38"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("E30").Value = "  -0.70%  "

$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("E32").Value = "  -0.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.63"
$ws.Range("E33").Value = "  +1.33%  "

$ws.Range("E34").Value = "  -2.92%  "

$ws.Range("E35").Value = "  +1.41%  "

$ws.Range("D36").Value = "0.0₃0820"
$ws.Range("E36").Value = "  -3.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.25"
$ws.Range("E37").Value = "  -3.45%  "

$ws.Range("E38").Value = "  -3.43%  "

$ws.Range("E39").Value = "  -0.98%  "

$ws.Range("E40").Value = "  -1.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.42"
$ws.Range("E41").Value = "  -0.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "433.40"
$ws.Range("E42").Value = "  -2.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.289"
$ws.Range("E43").Value = "  +0.82%  "

$ws.Range("E44").Value = "  +2.96%  "

$ws.Range("E45").Value = "  -0.79%  "

$ws.Range("D46").Value = "2.829.61"
$ws.Range("E46").Value = "  +0.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.20"
$ws.Range("E47").Value = "  -4.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.58"
$ws.Range("E48").Value = "  -2.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.07"
$ws.Range("E50").Value = "  -1.25%  "

$ws.Range("E51").Value = "  -1.64%  "

